$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rebuild the small table -------------------------------------------------
# The old sheet was a flat list of player names in column A (rows 1-8).
# The new sheet is a 1-row title plus a 5-column header row for an
# "elimination" tracking table.
$ws.Range("A1:A8").ClearContents()

$ws.Range("A1").Value = "Nom de la BDF"

$ws.Range("A2").Value = "Classement"
$ws.Range("B2").Value = "Joueur"
$ws.Range("C2").Value = "Heure "
$ws.Range("D2").Value = "Killer"
$ws.Range("E2").Value = "Points"

# Column A needs to be wide enough to show "Nom de la BDF" in full.
$ws.Columns.Item(1).ColumnWidth = 20.71

# --- Theme accent tweak (swap accent1 <-> accent5) ---------------------------
$tcs = $wb.Theme.ThemeColorScheme
$accent1 = $tcs.Colors(5).RGB
$accent5 = $tcs.Colors(9).RGB
$tcs.Colors(5).RGB = $accent5
$tcs.Colors(9).RGB = $accent1

# --- Footer carrying the sensitivity label text ------------------------------
$ws.PageSetup.CenterFooter = "`r&1#&`"Calibri`"&6&K626469 Public"

# --- Leave the selection where the author ended up ---------------------------
$ws.Range("C8").Select()
